$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column E (Weak_ban) and data/footnote column F
$ws.Range("E1").Value = "Weak_ban"

# Row -> (Weak_ban flag, footnote text) mapping.
# Rows without a footnote just get a 0 in column E.
$weakBanRows = @{
    7  = 'CT - "Stay Safe, Stay Home" order does not have any language that orders individuals to stay at home; 3/23/2020 safer-at-home order acted only to close nonessential businesses'
    9  = "GA - Only high-risk individuals were ordered to shelter in place"
    16 = "KY - Order on 3/25/2020 did not have any language ordering individuals to stay at home; 3/26/2020 order did not order residents to stay at home; Order never applied to the entire state but it expired officially on 6/29/2020"
    33 = "OK - Stay-at-home order from 4/1 only applied to those over 65; order was allowed to expire 5/15 per state's website"
    38 = "TX - Stay-at-home order issued 4/2, though it does not explicitly order individuals to stay at home"
}

for ($r = 2; $r -le 43; $r++) {
    $eCell = $ws.Cells.Item($r, 5)
    if ($weakBanRows.ContainsKey($r)) {
        $eCell.Value = 1
        $fCell = $ws.Cells.Item($r, 6)
        $fCell.Value = $weakBanRows[$r]
        $fCell.Font.Size = 10
    } else {
        $eCell.Value = 0
    }
    # Row 33 already carries a yellow highlight style (from A33:D33);
    # make sure the new E33 cell matches the rest of that row.
    if ($ws.Cells.Item($r, 4).Interior.ColorIndex -eq 6) {
        $eCell.Interior.Color = 65535
    }
}

$ws.Range("F35").Select() | Out-Null
